$wb = $excel.ActiveWorkbook

# --- Patients sheet: add "doctor_ID" column (N) ---
$patients = $wb.Worksheets.Item("Patients")
$doctors = $wb.Worksheets.Item("Doctors")

$patients.Cells.Item(1, 14).Value = "doctor_ID"
for ($r = 2; $r -le 20; $r++) {
    $patients.Cells.Item($r, 14).Value = $doctors.Cells.Item($r, 1).Text
}

# --- Departments sheet: insert a proper header row ---
$departments = $wb.Worksheets.Item("Departments")
$departments.Rows.Item(1).Insert()
$departments.Cells.Item(1, 1).Value = "department"

Write-Output "done"
